# "avancement sur la doc"
#
# Journal entry: a new work session was logged in row 31 of the
# "Activités" table (Tableau1) - the user started working on 2021-03-29
# at 11:00 and hasn't filled in "Fin" yet, so [Temps] keeps ticking via
# NOW(). Then the selection moved on to E32, the next blank "Activité"
# cell, ready for the next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date (col A, "Date") -> 2021-03-29
$ws.Range("A31").Value = 44284
# Début (col B, "Début") -> 11:00
$ws.Range("B31").Value = 0.45833333333333331

# Move the live selection to the next empty "Activité" cell.
$ws.Range("E32").Select()
